# Daily attendance processing - 2025-12-30 07:57:43
# Reorders the "Recorded By" (column G) entries on rows that include a
# "System"/"system" actor together with other recorders: the first
# recorder in the comma-separated list is moved to the end of the list.
# Rows with a single recorder, or rows whose list does not include
# "System"/"system" at all, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if ($value -eq "") { continue }

    $parts = $value -split ", "

    if ($parts.Count -le 1) { continue }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.Trim().ToLower() -eq "system") {
            $hasSystem = $true
        }
    }

    if (-not $hasSystem) { continue }

    $rotated = ($parts[1..($parts.Count - 1)] + $parts[0]) -join ", "
    $cell.Value = $rotated
}
